$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.194.10'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '1.861.30'
$ws.Range('E3').Value = '  -1.49%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '234.70'
$c.ClearFormats()
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('E6').Value = '  +0.02%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4658'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.67%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2825'
$c.ClearFormats()
$ws.Range('E8').Value = '  -0.93%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06524'
$c.ClearFormats()
$ws.Range('E9').Value = '  -1.42%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '21.27'
$c.ClearFormats()
$ws.Range('E10').Value = '  +3.21%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07844'
$c.ClearFormats()
$ws.Range('E11').Value = '  +0.60%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '97.17'
$c.ClearFormats()
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').Value = '1.865.38'
$ws.Range('E13').Value = '  -1.21%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '5.090'
$c.ClearFormats()
$ws.Range('E14').Value = '  -0.83%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.6714'
$c.ClearFormats()
$ws.Range('E15').Value = '  -0.83%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '278.43'
$c.ClearFormats()
$ws.Range('E16').Value = '  -2.58%  '
$ws.Range('D17').Value = '30.189.69'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('E18').Value = '  +0.02%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '5.509'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').Value = '2.111.57'
$ws.Range('E21').Value = '  -1.13%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.000007273'
$c.ClearFormats()
$ws.Range('E22').Value = '  -0.80%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.12%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.138'
$c.ClearFormats()
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.191'
$c.ClearFormats()
$ws.Range('E25').Value = '  -2.60%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '164.69'
$c.ClearFormats()
$ws.Range('E26').Value = '  -1.73%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '19.11'
$c.ClearFormats()
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('E28').Value = '  -4.28%  '
$ws.Range('E29').Value = '  -0.28%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.09677'
$c.ClearFormats()
$ws.Range('E30').Value = '  -0.96%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.416'
$c.ClearFormats()
$ws.Range('E31').Value = '  +0.14%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.473'
$c.ClearFormats()
$ws.Range('E32').Value = '  -1.23%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.074'
$c.ClearFormats()
$ws.Range('E33').Value = '  -2.36%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.04685'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.43%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.112'
$c.ClearFormats()
$ws.Range('E35').Value = '  +0.95%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.7047'
$c.ClearFormats()
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  +0.62%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.01850'
$c.ClearFormats()
$ws.Range('E38').Value = '  -1.74%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.529'
$c.ClearFormats()
$ws.Range('E39').Value = '  -0.08%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '6.225'
$c.ClearFormats()
$ws.Range('E40').Value = '  -7.38%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '73.12'
$c.ClearFormats()
$ws.Range('E41').Value = '  +0.46%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.941'
$c.ClearFormats()
$ws.Range('E42').Value = '  -2.18%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.8439'
$c.ClearFormats()
$ws.Range('E43').Value = '  -3.19%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '103.94'
$c.ClearFormats()
$ws.Range('E44').Value = '  -0.18%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.0000'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.03%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.4156'
$c.ClearFormats()
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '9.209'
$c.ClearFormats()
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '935.46'
$c.ClearFormats()
$ws.Range('E49').Value = '  -7.29%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '34.03'
$c.ClearFormats()
$ws.Range('E51').Value = '  -2.90%  '
